# Finished Week 13 logging
# Update Target Depth Data for Lions: OFF and DEF sheets, row 2 (H = Home totals)

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 374
$wsOff.Range("C2").Value = 279
$wsOff.Range("D2").Value = 72
$wsOff.Range("E2").Value = 32
$wsOff.Range("F2").Value = 6

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 354
$wsDef.Range("C2").Value = 263
$wsDef.Range("D2").Value = 84
$wsDef.Range("E2").Value = 42
